$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) contains numeric-looking text values (e.g. "313.83").
# Excel COM auto-converts such text to real numbers on assignment, which would
# change the cell type away from the original text/inlineStr representation.
# Temporarily force the Price column to Text format while writing the new values,
# then restore the original (default) style so no stray formatting is left behind.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"


# Row 2
$ws.Range("D2").Value = "27.316.27"
$ws.Range("E2").Value = "  +1.06%  "

# Row 3
$ws.Range("D3").Value = "1.856.47"
$ws.Range("E3").Value = "  +1.57%  "

# Row 4
$ws.Range("E4").Value = "  -0.65%  "

# Row 5
$ws.Range("D5").Value = "313.83"
$ws.Range("E5").Value = "  +0.82%  "

# Row 7
$ws.Range("D7").Value = "0.4618"
$ws.Range("E7").Value = "  -0.53%  "

# Row 8
$ws.Range("D8").Value = "0.3701"
$ws.Range("E8").Value = "  -1.17%  "

# Row 9
$ws.Range("D9").Value = "0.07315"
$ws.Range("E9").Value = "  +0.43%  "

# Row 10
$ws.Range("D10").Value = "0.8819"
$ws.Range("E10").Value = "  +1.90%  "

# Row 11
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").Value = "0.07816"
$ws.Range("E11").Value = "  -0.06%  "

# Row 12
$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").Value = "19.83"
$ws.Range("E12").Value = "  -0.83%  "

# Row 13
$ws.Range("D13").Value = "1.831.32"
$ws.Range("E13").Value = "  -0.06%  "

# Row 14
$ws.Range("D14").Value = "5.383"
$ws.Range("E14").Value = "  +0.53%  "

# Row 15
$ws.Range("D15").Value = "6.527"
$ws.Range("E15").Value = "  -0.32%  "

# Row 16
$ws.Range("D16").Value = "91.97"
$ws.Range("E16").Value = "  +0.07%  "

# Row 17
$ws.Range("E17").Value = "  -0.55%  "

# Row 18
$ws.Range("D18").Value = "0.000008867"
$ws.Range("E18").Value = "  +0.63%  "

# Row 20
$ws.Range("D20").Value = "14.80"
$ws.Range("E20").Value = "  +0.60%  "

# Row 21
$ws.Range("D21").Value = "27.341.23"
$ws.Range("E21").Value = "  +1.01%  "

# Row 22
$ws.Range("D22").Value = "5.116"
$ws.Range("E22").Value = "  -0.97%  "

# Row 23
$ws.Range("E23").Value = "  -0.99%  "

# Row 24
$ws.Range("D24").Value = "2.078.13"
$ws.Range("E24").Value = "  +0.04%  "

# Row 25
$ws.Range("D25").Value = "1.890"
$ws.Range("E25").Value = "  +2.47%  "

# Row 26
$ws.Range("D26").Value = "152.26"
$ws.Range("E26").Value = "  -0.79%  "

# Row 27
$ws.Range("E27").Value = "  +0.64%  "

# Row 28
$ws.Range("D28").Value = "2.066"
$ws.Range("E28").Value = "  -0.98%  "

# Row 29
$ws.Range("D29").Value = "5.114"

# Row 30
$ws.Range("D30").Value = "116.05"
$ws.Range("E30").Value = "  +0.37%  "

# Row 31
$ws.Range("D31").Value = "0.08858"
$ws.Range("E31").Value = "  -0.05%  "

# Row 32
$ws.Range("D32").Value = "0.7643"
$ws.Range("E32").Value = "  +4.55%  "

# Row 33
$ws.Range("D33").Value = "2.996"
$ws.Range("E33").Value = "  +0.92%  "

# Row 34
$ws.Range("E34").Value = "  +3.01%  "

# Row 35
$ws.Range("D35").Value = "4.490"
$ws.Range("E35").Value = "  +0.99%  "

# Row 36
$ws.Range("D36").Value = "2.617"
$ws.Range("E36").Value = "  +5.30%  "

# Row 37
$ws.Range("D37").Value = "0.01962"
$ws.Range("E37").Value = "  +0.84%  "

# Row 38
$ws.Range("E38").Value = "  -0.15%  "

# Row 39
$ws.Range("D39").Value = "2.989"
$ws.Range("E39").Value = "  +2.35%  "

# Row 40
$ws.Range("D40").Value = "0.05217"
$ws.Range("E40").Value = "  -0.40%  "

# Row 41
$ws.Range("D41").Value = "7.028"
$ws.Range("E41").Value = "  -4.07%  "

# Row 42
$ws.Range("D42").Value = "0.5152"
$ws.Range("E42").Value = "  -0.31%  "

# Row 43
$ws.Range("D43").Value = "0.1639"
$ws.Range("E43").Value = "  +0.53%  "

# Row 44
$ws.Range("D44").Value = "8.344"
$ws.Range("E44").Value = "  +1.60%  "

# Row 45
$ws.Range("D45").Value = "0.4834"
$ws.Range("E45").Value = "  +0.16%  "

# Row 46
$ws.Range("D46").Value = "10.27"
$ws.Range("E46").Value = "  +0.49%  "

# Row 47
$ws.Range("E47").Value = "  -0.64%  "

# Row 48
$ws.Range("D48").Value = "103.25"
$ws.Range("E48").Value = "  +0.37%  "

# Row 49
$ws.Range("D49").Value = "1.650"
$ws.Range("E49").Value = "  +1.59%  "

# Row 50
$ws.Range("D50").Value = "0.06225"

# Row 51
$ws.Range("D51").Value = "65.54"
$ws.Range("E51").Value = "  +1.74%  "

# Restore default styling on the Price column (removes the temporary text format)
$priceRange.Style = "Normal"
